# "catchup on work computer"
# The BuildingData sheet gets an AutoFilter applied on column A (Building),
# showing only two buildings: "12015 W Bluff Creek Drive" and "160 Spear".
# All the other data rows get hidden by the filter, and the sheet selection
# moves to one of the still-visible rows (A12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BuildingData")

# Apply the AutoFilter over the full used range, filtering column 1 (Building)
# down to the two values kept visible.
$ws.Range("A1:I68").AutoFilter(1, @("12015 W Bluff Creek Drive", "160 Spear"), 7)

# The filter range is recorded as the sheet-scoped hidden defined name
# "_xlnm._FilterDatabase", same as real Excel does when an AutoFilter is set.
$filterDb = $ws.Names.Add("_xlnm._FilterDatabase", "=BuildingData!`$A`$1:`$I`$68")
$filterDb.Visible = $false

# Selection moves onto one of the still-visible rows.
$ws.Range("A12").Select()
